# Applies the "Đồ án 3 5/7/2022" template edit:
#   1. Drop the stray _GoBack bookmark that sat before "Mã thẻ:".
#   2. Split "Chuyên ngành:" into "Đối tượng" + a relocated _GoBack
#      bookmark + ":".
#   3. Replace the trailing bookmarkEnd with a brand new table row
#      "Ngày hết hạn" / "${expiredDate}".

$d = $word.ActiveDocument
$tbl = $d.Tables(1)

# --- Step 1: remove the _GoBack bookmark anchored before "Mã thẻ:" ---
# The bookmark is not reachable through a plain text/range edit (Word
# keeps it pinned to the paragraph), so rebuild the first row: insert a
# fresh blank row with the same two cell values, then delete the
# original row that still carries the old bookmark.
$firstRow = $tbl.Rows(1)
$newFirstRow = $tbl.Rows.Add($firstRow)
$newFirstRow.Cells(1).Range.Text = 'Mã thẻ:'
$newFirstRow.Cells(2).Range.Text = '${idStudent}'
$tbl.Rows(2).Delete()

# --- Step 2: "Chuyên ngành:" -> "Đối tượng" + bookmark + ":" ---
$rng = $d.Content
$rng.Find.Execute('Chuyên ngành:', $true, $false, $false, $false, $false, $true, 1, $false, 'Đối tượng:', 2) | Out-Null

$findLabel = $d.Content
$findLabel.Find.Execute('Đối tượng', $true, $false, $false, $false, $false, $true, 1, $false, '', 0) | Out-Null
$splitPos = $findLabel.End
$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add('_GoBack', $bmRange)

# --- Step 3: replace the trailing bookmarkEnd with a new table row ---
$lastRow = $tbl.Rows($tbl.Rows.Count)
$expRow = $tbl.Rows.Add()
$expRow.Cells(1).Range.Text = 'Ngày hết hạn'
$expRow.Cells(2).Range.Text = '${expiredDate}'
